$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.829945333333333
$ws.Range("H2").Value = 5.489835999999999
$ws.Range("I2").Value = 0.4190796720210465
$ws.Range("J2").Value = 0.4190796720210465
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 307.6317940161057
$ws.Range("R2").Value = 2768.686146144952
$ws.Range("S2").Value = 0.1250610121947294
$ws.Range("T2").Value = 0.1250610121947294

# Row 3
$ws.Range("G3").Value = 1.829945333333333
$ws.Range("H3").Value = 5.489835999999999
$ws.Range("I3").Value = 0.4190796720210465
$ws.Range("J3").Value = 0.4190796720210465
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 298.2925002624502
$ws.Range("R3").Value = 2684.632502362052
$ws.Range("S3").Value = 0.1212643255299079
$ws.Range("T3").Value = 0.1212643255299079

# Row 4
$ws.Range("G4").Value = 1.829945333333333
$ws.Range("H4").Value = 5.489835999999999
$ws.Range("I4").Value = 0.4190796720210465
$ws.Range("J4").Value = 0.4190796720210465
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 303.7590855870266
$ws.Range("R4").Value = 2733.83177028324
$ws.Range("S4").Value = 0.1234866468479203
$ws.Range("T4").Value = 0.1234866468479203

# Row 5
$ws.Range("G5").Value = 1.829945333333333
$ws.Range("H5").Value = 5.489835999999999
$ws.Range("I5").Value = 0.4190796720210465
$ws.Range("J5").Value = 0.4190796720210465
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 121.1913034351898
$ws.Range("R5").Value = 1090.721730916708
$ws.Range("S5").Value = 0.04926768744848892
$ws.Range("T5").Value = 0.04926768744848892

# Row 6
$ws.Range("I6").Value = 0.2833335737960661
$ws.Range("J6").Value = 0.2833335737960661
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 207.9853102669726
$ws.Range("R6").Value = 1871.867792402754
$ws.Range("S6").Value = 0.08455190240271679
$ws.Range("T6").Value = 0.0845519024027168

# Row 7
$ws.Range("I7").Value = 0.2833335737960661
$ws.Range("J7").Value = 0.2833335737960661
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.08198501864970639
$ws.Range("T7").Value = 0.08198501864970639

# Row 8
$ws.Range("I8").Value = 0.2833335737960661
$ws.Range("J8").Value = 0.2833335737960661
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 205.36702932247
$ws.Range("R8").Value = 1848.30326390223
$ws.Range("S8").Value = 0.08348749725507289
$ws.Range("T8").Value = 0.08348749725507289

# Row 9
$ws.Range("I9").Value = 0.2833335737960661
$ws.Range("J9").Value = 0.2833335737960661
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 81.93564949046566
$ws.Range("R9").Value = 737.420845414191
$ws.Range("S9").Value = 0.03330915548857
$ws.Range("T9").Value = 0.03330915548857

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1530633333333333
$ws.Range("H10").Value = 0.45919
$ws.Range("I10").Value = 0.03505335944376924
$ws.Range("J10").Value = 0.03505335944376924
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 25.73145053773111
$ws.Range("R10").Value = 231.58305483958
$ws.Range("S10").Value = 0.01046056133365328
$ws.Range("T10").Value = 0.01046056133365328

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1530633333333333
$ws.Range("H11").Value = 0.45919
$ws.Range("I11").Value = 0.03505335944376924
$ws.Range("J11").Value = 0.03505335944376924
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 24.95027778525889
$ws.Range("R11").Value = 224.55250006733
$ws.Range("S11").Value = 0.01014299254842556
$ws.Range("T11").Value = 0.01014299254842556

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1530633333333333
$ws.Range("H12").Value = 0.45919
$ws.Range("I12").Value = 0.03505335944376924
$ws.Range("J12").Value = 0.03505335944376924
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 25.40752301356666
$ws.Range("R12").Value = 228.6677071221
$ws.Range("S12").Value = 0.01032887564694036
$ws.Range("T12").Value = 0.01032887564694036

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1530633333333333
$ws.Range("H13").Value = 0.45919
$ws.Range("I13").Value = 0.03505335944376924
$ws.Range("J13").Value = 0.03505335944376924
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 10.13688471284111
$ws.Range("R13").Value = 91.23196241557
$ws.Range("S13").Value = 0.004120929914750027
$ws.Range("T13").Value = 0.004120929914750027

# Row 14
$ws.Range("G14").Value = 1.146373333333333
$ws.Range("H14").Value = 3.43912
$ws.Range("I14").Value = 0.2625333947391181
$ws.Range("J14").Value = 0.2625333947391181
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 192.7166231262044
$ws.Range("R14").Value = 1734.44960813584
$ws.Range("S14").Value = 0.07834474987215242
$ws.Range("T14").Value = 0.07834474987215244

# Row 15
$ws.Range("G15").Value = 1.146373333333333
$ws.Range("H15").Value = 3.43912
$ws.Range("I15").Value = 0.2625333947391181
$ws.Range("J15").Value = 0.2625333947391181
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 186.8660017353155
$ws.Range("R15").Value = 1681.79401561784
$ws.Range("S15").Value = 0.07596630704749956
$ws.Range("T15").Value = 0.07596630704749956

# Row 16
$ws.Range("G16").Value = 1.146373333333333
$ws.Range("H16").Value = 3.43912
$ws.Range("I16").Value = 0.2625333947391181
$ws.Range("J16").Value = 0.2625333947391181
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 190.2905562978666
$ws.Range("R16").Value = 1712.6150066808
$ws.Range("S16").Value = 0.07735848519110945
$ws.Range("T16").Value = 0.07735848519110945

# Row 17
$ws.Range("G17").Value = 1.146373333333333
$ws.Range("H17").Value = 3.43912
$ws.Range("I17").Value = 0.2625333947391181
$ws.Range("J17").Value = 0.2625333947391181
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 75.92056219348444
$ws.Range("R17").Value = 683.28505974136
$ws.Range("S17").Value = 0.0308638526283567
$ws.Range("T17").Value = 0.0308638526283567
